$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 1; $i -le 23; $i++) {
    $row = 18 + $i
    $ws.Range("A$row").Value = "Pre-or-Post-Chorus $i"
}
